$wb = $excel.ActiveWorkbook

# --- Step 1: reserve a sheetId slot so the upcoming duplicate lands on sheetId=3 ---
# (a throwaway sheet is inserted before the original; it is deleted again at the end)
$placeholder = $wb.Worksheets.Add()

# --- Step 2: duplicate the original Kanban sheet, placing the copy right after it ---
$original = $wb.Worksheets.Item(2)
$original.Copy($null, $original)

# --- Step 3: rename the two sheets ---
$startSheet = $wb.Worksheets.Item(2)
$startSheet.Name = "Начало"
$endSheet = $wb.Worksheets.Item(3)
$endSheet.Name = "Конец"

# --- Step 4: on "Конец", move every task card from the To-Do column (B) to the
#             Done column (D), keeping each row's original card formatting ---
$endSheet.Range("B3:B14").Copy($endSheet.Range("D3:D14"))
$endSheet.Range("B3:B14").ClearContents()

# widen the Done column now that it holds the task text
$endSheet.Columns.Item(4).ColumnWidth = 33.14

# reflect the new selection on the finished board
$endSheet.Range("I5").Select()

# --- Step 5: drop the throwaway placeholder sheet ---
$toRemove = $wb.Worksheets.Item(1)
$toRemove.Delete()

# --- Step 6: make "Конец" the active/visible tab ---
$final = $wb.Worksheets.Item(2)
$final.Activate()
